$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Snapshot every existing hyperlink (row/column/address) before we disturb
#    the sheet. We'll recreate them after the insert so their IDs regenerate
#    in the same relative order (rId1..rId37), then append the new one last
#    (rId38).
# ---------------------------------------------------------------------------
$existingLinks = @()
foreach ($hl in $ws.Hyperlinks) {
  $rng = $hl.Range
  $existingLinks += [PSCustomObject]@{
    Row     = $rng.Row
    Col     = $rng.Column
    Address = $hl.Address
  }
}

# Remove all hyperlinks so we can re-add them (this engine's Hyperlinks.Delete
# operates sheet-wide regardless of which Range/collection it was fetched
# from, so do it once up front).
$ws.Hyperlinks.Delete()

# ---------------------------------------------------------------------------
# 2. Insert the new row for "Joint simplified divorce - Cook County" above
#    the existing "Living will" row (row 23), pushing everything else down.
# ---------------------------------------------------------------------------
$ws.Range("A23").EntireRow.Insert()

# Write URL before name so the shared-string table appends the url text
# ahead of the display name (matches the authored workbook ordering).
$ws.Range("B23").Value = "https://www.illinoislegalaid.org/legal-information/joint-simplified-divorce-cook-county"
$ws.Range("A23").Value = "Joint simplified divorce - Cook County"

# ---------------------------------------------------------------------------
# 3. Recreate the original hyperlinks, shifting any row that was at/after 23
#    down by one to follow its data.
# ---------------------------------------------------------------------------
foreach ($e in $existingLinks) {
  $newRow = $e.Row
  if ($newRow -ge 23) {
    $newRow = $newRow + 1
  }
  $target = $ws.Cells.Item($newRow, $e.Col)
  $ws.Hyperlinks.Add($target, $e.Address) | Out-Null
  # Adding a hyperlink can reassign a fresh cell style; put it back on the
  # shared "Hyperlink" style so we don't fork a duplicate style entry per
  # cell (matches the original workbook, where every linked B cell uses the
  # same style index).
  $target.Style = "Hyperlink"
}

# New hyperlink for the inserted row, added last so it becomes rId38.
$ws.Hyperlinks.Add($ws.Range("B23"), "https://www.illinoislegalaid.org/legal-information/joint-simplified-divorce-cook-county") | Out-Null
$ws.Range("B23").Style = "Hyperlink"

# Leave the selection on the newly-added row, like the author would have
# after typing the new entry in.
$ws.Range("A23").Select() | Out-Null
